$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet in front of "ODI Batting" (which is
#    the active sheet, so Worksheets.Add() inserts right before it).
#    NOTE: worksheet references are index-bound, so any sheet reference
#    fetched *before* this insertion would silently start pointing at the
#    newly inserted sheet afterwards. Always re-fetch sheets by name once
#    the sheet collection has been restructured.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Header row (bold, thin border, centered / top aligned - mirrors the header
# style already used on the other sheets).
$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous (thin)

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row - keep the ID as text (matches the source data, which stores it as
# a plain string rather than a number).
$idCell = $playerInfo.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "6544"
$idCell.Style = "Normal"

$playerInfo.Range("B2").Value = "Kuldeep Rampal Sen"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet: MATCH_CARD_LINK column becomes MATCH_CODE, and the
#    stored value becomes the bare match code instead of the full URL.
#    (Fetched fresh by name now that the sheet list has been restructured.)
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCode = $battingSheet.Range("D2")
$battingCode.NumberFormat = "@"
$battingCode.Value = "4679"
$battingCode.Style = "Normal"

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet: same MATCH_CARD_LINK -> MATCH_CODE change.
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCode = $bowlingSheet.Range("B2")
$bowlingCode.NumberFormat = "@"
$bowlingCode.Value = "4679"
$bowlingCode.Style = "Normal"
